$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D that look like plain numbers need to be forced to
# text so Excel keeps them as strings (matching the source inlineStr cells)
# instead of auto-converting them to numeric values. The NumberFormat is
# reset back to the default afterwards so no stray cell formatting remains.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "26.610.95"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").Value = "1.686.59"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  -0.13%  "
Set-TextValue $ws.Range("D5") "217.79"
$ws.Range("E5").Value = "  +4.04%  "
Set-TextValue $ws.Range("D6") "0.5330"
$ws.Range("E6").Value = "  +2.96%  "
$ws.Range("E7").Value = "  -0.10%  "
Set-TextValue $ws.Range("D8") "0.2684"
$ws.Range("E8").Value = "  +5.08%  "
Set-TextValue $ws.Range("D9") "0.06435"
$ws.Range("E9").Value = "  +3.31%  "
Set-TextValue $ws.Range("D10") "21.55"
$ws.Range("E10").Value = "  +6.78%  "
Set-TextValue $ws.Range("D11") "0.07803"
$ws.Range("E11").Value = "  +3.27%  "
$ws.Range("D12").Value = "1.693.17"
$ws.Range("E12").Value = "  +3.70%  "
Set-TextValue $ws.Range("D13") "4.511"
$ws.Range("E13").Value = "  +3.71%  "
Set-TextValue $ws.Range("D14") "0.5649"
$ws.Range("E14").Value = "  +4.66%  "
$ws.Range("D15").Value = "0.0₅8473"
$ws.Range("E15").Value = "  +7.54%  "
Set-TextValue $ws.Range("D16") "66.46"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "26.654.23"
$ws.Range("E17").Value = "  +2.87%  "
Set-TextValue $ws.Range("D18") "1.002"
$ws.Range("E18").Value = "  -0.14%  "
Set-TextValue $ws.Range("D19") "4.810"
$ws.Range("E19").Value = "  +4.28%  "
Set-TextValue $ws.Range("D20") "196.08"
$ws.Range("E20").Value = "  +6.71%  "
Set-TextValue $ws.Range("D21") "10.43"
$ws.Range("E21").Value = "  +4.72%  "
Set-TextValue $ws.Range("D22") "6.391"
$ws.Range("E22").Value = "  +5.50%  "
$ws.Range("E23").Value = "  -0.14%  "
Set-TextValue $ws.Range("D24") "143.69"
$ws.Range("E24").Value = "  -0.79%  "
$ws.Range("E25").Value = "  +7.70%  "
Set-TextValue $ws.Range("D26") "7.486"
$ws.Range("E26").Value = "  +2.42%  "
Set-TextValue $ws.Range("D27") "16.27"
$ws.Range("E27").Value = "  +5.46%  "
Set-TextValue $ws.Range("D29") "0.06196"
$ws.Range("E29").Value = "  +4.71%  "
Set-TextValue $ws.Range("D30") "1.282"
$ws.Range("E30").Value = "  +3.28%  "
$ws.Range("E31").Value = "  +8.55%  "
Set-TextValue $ws.Range("D32") "3.473"
$ws.Range("E32").Value = "  +4.07%  "
Set-TextValue $ws.Range("D33") "1.704"
$ws.Range("E33").Value = "  +6.51%  "
Set-TextValue $ws.Range("D34") "1.015"
$ws.Range("E34").Value = "  +5.03%  "
$ws.Range("E35").Value = "  +2.33%  "
Set-TextValue $ws.Range("D36") "2.420"
$ws.Range("E36").Value = "  +1.56%  "
Set-TextValue $ws.Range("D37") "0.5726"
$ws.Range("E37").Value = "  -1.33%  "
Set-TextValue $ws.Range("D38") "0.01651"
$ws.Range("E38").Value = "  +3.97%  "
Set-TextValue $ws.Range("D39") "6.006"
$ws.Range("E39").Value = "  +6.36%  "
$ws.Range("D40").Value = "1.077.24"
$ws.Range("E40").Value = "  +4.55%  "
Set-TextValue $ws.Range("D41") "0.8661"
$ws.Range("E41").Value = "  +3.36%  "
$ws.Range("E42").Value = "  -0.01%  "
Set-TextValue $ws.Range("D43") "100.52"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").Value = "1.837.39"
$ws.Range("E44").Value = "  +3.06%  "
$ws.Range("D45").Value = "0.0₈109"
$ws.Range("E45").Value = "  +3.40%  "
Set-TextValue $ws.Range("D46") "57.51"
$ws.Range("E46").Value = "  +6.35%  "
Set-TextValue $ws.Range("D47") "8.168"
$ws.Range("E47").Value = "  +2.60%  "
Set-TextValue $ws.Range("D48") "1.006"
$ws.Range("E48").Value = "  +0.54%  "
Set-TextValue $ws.Range("D49") "0.05219"
$ws.Range("E49").Value = "  +0.80%  "
Set-TextValue $ws.Range("D50") "6.093"
$ws.Range("E50").Value = "  +5.62%  "
Set-TextValue $ws.Range("D51") "0.4240"
$ws.Range("E51").Value = "  +0.24%  "
